# "cambio en el link de la pagina" — update the URL text shown/stored in C10.
# The cell keeps its existing hyperlink relationship (still pointing at the
# github.com repo URL); only the displayed text/value changes to the new
# GitHub Pages URL.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "https://contesl.github.io/C24172G11/"
